# Update the price list date and the two unit prices on "Hoja1".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header date serial (A1): 24/04/2024 -> 20/05/2024
$ws.Range("A1").Value = 45432

# "TOMA para Canilla de GOMA de 1/2" price
$ws.Range("D30").Value = 202.54

# "TOMA para Canilla de GOMA de 3/4" price
$ws.Range("D31").Value = 230
